# Corrected test cases for .Drop
#
# The ".Drop" example table had a bogus extra row testing [1,2] (which
# doesn't belong here - that case belongs to other ops) and the existing
# [1] row reported the wrong output ("[], 0" instead of "[], 1" - dropping
# a 1-element stack leaves the dropped value 1, not 0).
#
# Fix: remove the erroneous "[1,2]" row from the .Drop table (this shifts
# the following ".Calculate" block up by one row) and correct the output
# of the "[1]" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the spurious "[1,2] -> [1], 2 -> [1]" row from the .Drop test
# table; everything below (the .Calculate section) shifts up one row.
$ws.Rows("12:12").Delete()

# Correct the .Drop output for input stack [1]: dropping its single
# element leaves an empty stack and the dropped value is 1.
$ws.Range("B11").Value = "[], 1"

# Move the active selection, matching the saved view state.
$null = $ws.Range("N3").Select()
